# Apply quarterly database update + read_price algorithm shift.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Shift the quarter-label headers (row 8 and row 24, columns E:N)
#    one quarter forward: the oldest quarter ("فصل دوم منتهی به
#    1399/06") drops off and a new quarter ("فصل چهارم منتهی به
#    1401/12") is appended at the end.
# ---------------------------------------------------------------------
$quarterLabels = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)

$headerRows = @(8, 24)
foreach ($r in $headerRows) {
    for ($i = 0; $i -lt $quarterLabels.Length; $i++) {
        $col = 5 + $i   # column E = 5 ... N = 14
        $ws.Cells.Item($r, $col).Value = $quarterLabels[$i]
    }
}

# ---------------------------------------------------------------------
# 2) Shift the data rows left by one quarter (drop oldest value in E,
#    shift remaining values left, and place the newly reported quarter
#    value in N).
# ---------------------------------------------------------------------
$rowData = @{
    10 = @(845, 1178, 646, 755, 1115, 1451, 677, 1486, 1218, 3009)
    13 = @(3782, 1935, 826, 241, 15785, 46750, 8692, 21696, 31931, -3419)
    14 = @(1766, 2905, 2800, 3118, 2959, 6924, 8540, 3277, 7831, 15273)
    15 = @(188, 508, 165, 211, 188, 252, 176, 238, -118, 386)
    16 = @(3290, 1758, 1509, 1547, 1528, 1557, 1623, 1967, 2361, 2680)
    17 = @(18580, 39828, 40275, 25036, 35310, 47620, 67982, 50371, 18643, 119316)
    19 = @(12685, 20132, 17289, 53727, 27933, 99970, 37609, 127681, -14508, 66594)
    20 = @(41136, 68244, 63510, 84635, 84818, 204524, 125299, 206716, 47358, 203839)
    26 = @(229, 187, 209, 208, 207, 234, 232, 235, 237, 241)
    27 = @(153, 197, 174, 174, 175, 147, 143, 145, 145, 143)
}

foreach ($r in $rowData.Keys) {
    $values = $rowData[$r]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 5 + $i   # column E = 5 ... N = 14
        $ws.Cells.Item($r, $col).Value = $values[$i]
    }
}

$wb.Save()
